$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update subcategory (H column) text for rows with revised labels
$ws.Range("H3").Value = "line graph(s)"
$ws.Range("H4").Value = "line graph(s)"
$ws.Range("H5").Value = "line graph(s)"
$ws.Range("H6").Value = "line graph(s)"
$ws.Range("H7").Value = "scatter plot(s)"
$ws.Range("H12").Value = "line graph(s)"
$ws.Range("H13").Value = "line graph(s)"
$ws.Range("H14").Value = "line graph(s)"
$ws.Range("H17").Value = "photo(s)"
$ws.Range("H25").Value = "photo(s)"
$ws.Range("H26").Value = "photo(s)"
$ws.Range("H27").Value = "photo(s)"
$ws.Range("H28").Value = "photo(s)"
$ws.Range("H30").Value = "photo(s)"
$ws.Range("H31").Value = "photo(s)"
$ws.Range("H33").Value = "line graph(s)"
$ws.Range("H36").Value = "mixed statistical plot (more than 1 statistical plot and type)"
$ws.Range("H37").Value = "line graph(s)"
$ws.Range("H39").Value = "data display"
$ws.Range("H41").Value = "photo(s)"
$ws.Range("H43").Value = "data collection, data analysis, data gathering diagram"
$ws.Range("H44").Value = "photo(s)"
$ws.Range("H45").Value = "photo(s)"
$ws.Range("H46").Value = "photo(s)"
$ws.Range("H51").Value = "drawing(s)"
$ws.Range("H57").Value = "line graph(s)"
$ws.Range("H58").Value = "line graph(s)"
$ws.Range("H59").Value = "line graph(s)"
$ws.Range("H60").Value = "line graph(s)"
$ws.Range("H61").Value = "line graph(s)"
$ws.Range("H62").Value = "line graph(s)"
$ws.Range("H69").Value = "line graph(s)"
$ws.Range("H70").Value = "scatter plot(s)"
$ws.Range("H76").Value = "line graph(s)"
$ws.Range("H77").Value = "line graph(s)"
$ws.Range("H80").Value = "data display"
$ws.Range("H81").Value = "line graph(s)"
$ws.Range("H85").Value = "drawing(s)"

# Remove the is_viewed column (I) entirely, shifting dimension back to A1:H93
$ws.Range("I1").EntireColumn.Delete()

Write-Output "done"
